# Insert a new weekly price record for Coco (Mercado Mayorista Lo Valledor de
# Santiago) as row 13, pushing all the existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13; everything from row 13 down moves to
# row 14 and beyond (dimension grows from T54 to T55 automatically).
$ws.Rows("13:13").Insert()

# Fill the newly inserted row 13 with the new record's data.
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44459
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100108
$ws.Range("H13").Value = "Tropicales y subtropicales"
$ws.Range("I13").Value = 100108007
$ws.Range("J13").Value = "Coco"
$ws.Range("K13").Value = "Sin especificar"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 130
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("Q13").Value = "$/malla 20 unidades"
$ws.Range("R13").Value = "Perú"
$ws.Range("S13").Value = 1025
$ws.Range("T13").Value = 20
